# Update the descriptive statistics for "Future time perspective" (futi_mean),
# row 31, after the R environment activation procedure / working directory
# change caused the underlying computation to produce slightly different
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("D31").Value = 2.6451137884872802
$ws.Range("E31").Value = 0.69439933930207298
$ws.Range("G31").Value = 4.9000000000000004
$ws.Range("H31").Value = 2.1749999999999998
$ws.Range("I31").Value = 2.65
$ws.Range("J31").Value = 3.1
$ws.Range("K31").Value = 0.34315607081937599
$ws.Range("L31").Value = 0.13211102018389001
